$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Row 11 new values
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = 0.8994

# Row 21 new values
$ws.Range("D21").Value = 86
$ws.Range("E21").Value = 88
$ws.Range("F21").Value = 841
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 149.0453

# Row 31 new values
$ws.Range("D31").Value = 316
$ws.Range("E31").Value = 318
$ws.Range("F31").Value = 2912
$ws.Range("G31").Value = 12
$ws.Range("H31").Value = 764.6784

# Update active selection on the Results sheet
$ws.Activate()
$ws.Range("H32").Select()
